$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = "maa://22742 (91.41), *maa://20791 (63.01)"
$ws.Range("T3").Value = "maa://24617 (89.57), **maa://20790 (43.48), ***maa://37170 (17.19), maa://45854 (100.0)"
$ws.Range("X3").Value = "maa://27396 (84.28), maa://27484 (96.46), maa://27480 (82.86)"
$ws.Range("AB3").Value = "maa://24390 (94.2)"
$ws.Range("X4").Value = "**maa://32495 (48.51), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (86.67)"
$ws.Range("D6").Value = "maa://42407 (95.35)"
$ws.Range("AF6").Value = "*maa://33152 (59.57), ***maa://22770 (26.09)"
$ws.Range("AF7").Value = "*maa://26191 (69.05), *maa://36671 (68.0), *maa://42530 (62.5), maa://45272 (100.0)"
$ws.Range("A8").Value = "更新日期：2025.01.31 13:17:43"
$ws.Range("D9").Value = "maa://22765 (92.47), *maa://21915 (70.37)"
$ws.Range("AB9").Value = "maa://28711 (86.61), ***maa://22740 (5.77), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (95.24)"
$ws.Range("D10").Value = "***maa://25695 (18.82), **maa://32237 (41.3), ***maa://34206 (20.0), ***maa://39951 (14.89), ***maa://39243 (28.57), *maa://45271 (52.38)"
$ws.Range("T11").Value = "maa://22747 (92.9), maa://22501 (97.59), *maa://45521 (71.43)"
$ws.Range("X11").Value = "maa://36713 (97.88)"
$ws.Range("AF11").Value = "maa://31203 (96.0)"
$ws.Range("D13").Value = "maa://24999 (91.9), maa://36673 (92.96), maa://25001 (85.51)"
$ws.Range("H13").Value = "*maa://21248 (73.39), **maa://22728 (47.73)"
$ws.Range("P13").Value = "maa://22676 (92.62), *maa://22583 (74.63), *maa://22500 (58.7)"
$ws.Range("T14").Value = "maa://22521 (94.06), maa://42751 (100.0)"
$ws.Range("D15").Value = "*maa://22743 (77.45), maa://22734 (84.03), *maa://30808 (64.18), **maa://36048 (42.11), maa://45058 (100.0)"
$ws.Range("P15").Value = "maa://24762 (90.62), *maa://22727 (70.0)"
$ws.Range("D16").Value = "maa://21441 (96.4), maa://36679 (93.88), maa://37650 (97.06)"
$ws.Range("T16").Value = "maa://22729 (94.84), *maa://28648 (69.23), maa://36674 (80.43)"
$ws.Range("H17").Value = "maa://22430 (88.66), maa://39599 (85.42)"
$ws.Range("P17").Value = "maa://23890 (81.19), *maa://24940 (67.86)"
$ws.Range("D18").Value = "maa://24570 (97.25)"
$ws.Range("L18").Value = "maa://22466 (89.54), *maa://22732 (51.16)"
$ws.Range("AF18").Value = "*maa://24313 (58.64), **maa://29784 (44.44)"
$ws.Range("AB19").Value = "*maa://30709 (64.79), *maa://36668 (57.5)"
$ws.Range("AF19").Value = "*maa://21663 (63.38)"
$ws.Range("L20").Value = "maa://41331 (85.29)"
$ws.Range("AB21").Value = "maa://21443 (80.76), ***maa://23820 (29.31)"
$ws.Range("AF21").Value = "maa://22524 (94.39), *maa://22432 (77.61)"
$ws.Range("X22").Value = "maa://21282 (98.58), *maa://37649 (67.86)"
$ws.Range("AF22").Value = "maa://29658 (93.62)"
$ws.Range("L23").Value = "maa://39756 (95.25), maa://39875 (94.12)"
$ws.Range("P23").Value = "maa://30587 (91.79), *maa://29748 (75.78), ***maa://29785 (16.42), *maa://37566 (74.29)"
$ws.Range("X24").Value = "maa://29988 (85.08), maa://23504 (93.19), **maa://22892 (40.14), *maa://25141 (76.74), *maa://36663 (78.08), ***maa://22815 (23.08)"
$ws.Range("AF25").Value = "maa://20108 (96.3), maa://24621 (96.75), maa://36676 (96.97), maa://22771 (85.71), *maa://37772 (66.67)"
$ws.Range("D26").Value = "maa://41802 (93.33)"
$ws.Range("AB26").Value = "maa://42235 (94.57)"
$ws.Range("D28").Value = "maa://24465 (90.99), maa://25725 (83.72)"
$ws.Range("X28").Value = "maa://39929 (90.44), maa://41749 (92.0), ***maa://39723 (13.89)"
$ws.Range("L29").Value = "maa://28432 (93.21), *maa://28440 (79.05), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("P29").Value = "*maa://23168 (57.63), *maa://30050 (51.61)"
$ws.Range("X30").Value = "maa://39477 (88.24)"
$ws.Range("L31").Value = "maa://35926 (93.55), maa://36258 (85.32), *maa://43904 (72.73)"
$ws.Range("T32").Value = "maa://42859 (96.12), maa://41108 (88.0), maa://41238 (97.0), maa://45523 (100.0)"
$ws.Range("L35").Value = "maa://41296 (95.97)"
$ws.Range("T39").Value = "maa://45788 (82.35), maa://45790 (88.89)"
$ws.Range("P41").Value = "**maa://35616 (38.24), maa://43177 (88.89)"
$ws.Range("T44").Value = "maa://39366 (88.24)"
$ws.Range("H47").Value = "maa://27410 (96.38), maa://29661 (97.24), maa://28038 (84.62)"
$ws.Range("H53").Value = "maa://32534 (93.77), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (92.14)"
$ws.Range("H59").Value = "maa://27746 (82.73), maa://31270 (95.2)"
$ws.Range("H60").Value = "*maa://40438 (67.86)"
$ws.Range("H62").Value = "maa://42981 (97.3), maa://43903 (100.0)"
